$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "hören",
    "siegen",
    "fällen",
    "stammen",
    "jagen",
    "trennen",
    "altern",
    "schmecken",
    "seufzen",
    "ärgern",
    "machen",
    "treiben",
    "schreiten",
    "spüren",
    "wachsen",
    "weichen",
    "zünden",
    "dringen",
    "heben",
    "bluten",
    "streichen",
    "wüten",
    "werfen",
    "scheitern",
    "lesen",
    "kichern",
    "schenken",
    "enden",
    "quälen",
    "tollen",
    "fließen",
    "planen",
    "schrecken",
    "decken",
    "arten",
    "irren",
    "folgen",
    "kümmern",
    "liefern",
    "sorgen",
    "tropfen",
    "räumen",
    "drehen",
    "greifen",
    "hauen",
    "sperren",
    "werden",
    "lockern",
    "töten",
    "bauen",
    "platzen",
    "schlucken",
    "sterben",
    "wirken",
    "schulden",
    "klingen",
    "führen",
    "trauen",
    "rasen",
    "zielen",
    "stecken",
    "klettern",
    "boxen",
    "wundern",
    "fischen",
    "lügen",
    "schwingen",
    "helfen",
    "warnen",
    "reizen",
    "freuen",
    "fangen",
    "flehen",
    "kehren",
    "suchen",
    "ändern",
    "bergen",
    "äußern",
    "jubeln",
    "biegen",
    "graben",
    "pfeifen",
    "scheinen",
    "münzen",
    "malen",
    "schwächen",
    "backen",
    "feiern",
    "spinnen",
    "erben",
    "formen",
    "stehlen",
    "flüchten",
    "achten",
    "wehtun",
    "sinken",
    "mauern",
    "spielen",
    "sichern",
    "fahren",
    "zeigen",
    "mögen",
    "sprengen",
    "heilen",
    "bellen",
    "ehren",
    "bitten",
    "gelten",
    "gründen",
    "liegen",
    "heulen",
    "loben",
    "filmen",
    "dienen",
    "wenden",
    "rufen",
    "saufen",
    "grüßen",
    "knarren",
    "schwören",
    "geben",
    "runden",
    "brauchen",
    "betteln",
    "pflanzen",
    "fallen",
    "kosten",
    "zögern"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $words[$i]
}

